$wb = $excel.ActiveWorkbook

# A cell elsewhere in the workbook already holds the literal text
# "2012-04-24" as a shared string (the "date" column of the land sheet).
# Copying it with PasteSpecial(xlPasteValues) lets us reuse that exact
# string wherever we need the same text without Excel re-interpreting the
# typed text "2012-04-24" as a real date serial number.
$dateSource = $wb.Worksheets.Item(1).Range("K2")

# ---------------------------------------------------------------------------
# Sheet 3 ("現金" / cash): turn the buggy "header row == data row" into real
# field-name headers, and append the per-record metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other "normal" sheets already carry.
# ---------------------------------------------------------------------------
$wsCash = $wb.Worksheets.Item(3)

# Give the new header cells (E1:K1) the same style as the existing header
# cells (B1:D1) before filling them in.
$wsCash.Range("B1:D1").Copy()
$wsCash.Range("E1:K1").PasteSpecial(-4122)

$wsCash.Range("B1").Value = "currency"
$wsCash.Range("C1").Value = "owner"
$wsCash.Range("D1").Value = "total"
$wsCash.Range("E1").Value = "property_category"
$wsCash.Range("F1").Value = "category"
$wsCash.Range("G1").Value = "date"
$wsCash.Range("H1").Value = "legislator_name"
$wsCash.Range("I1").Value = "legislator_id"
$wsCash.Range("J1").Value = "source_file"
$wsCash.Range("K1").Value = "index"

# Give the new data cells (E2:K2) the same style as the existing data cells
# (B2:D2) before filling them in.
$wsCash.Range("B2:D2").Copy()
$wsCash.Range("E2:K2").PasteSpecial(-4122)

$wsCash.Range("E2").Value = "cash"
$wsCash.Range("F2").Value = "normal"
$dateSource.Copy()
$wsCash.Range("G2").PasteSpecial(-4163)
$wsCash.Range("H2").Value = "段宜康"
$wsCash.Range("I2").Value = 917
$wsCash.Range("J2").Value = "tmp25ce1"
$wsCash.Range("K2").Value = 45

# ---------------------------------------------------------------------------
# Sheet 4 ("存款" / deposit): same story — real headers instead of the
# copied-data-row bug, plus the metadata columns G:M on every data row.
# ---------------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item(4)

$wsDeposit.Range("B1:F1").Copy()
$wsDeposit.Range("G1:M1").PasteSpecial(-4122)

$wsDeposit.Range("B1").Value = "bank"
$wsDeposit.Range("C1").Value = "deposit_type"
$wsDeposit.Range("D1").Value = "currency"
$wsDeposit.Range("E1").Value = "owner"
$wsDeposit.Range("F1").Value = "total"
$wsDeposit.Range("G1").Value = "property_category"
$wsDeposit.Range("H1").Value = "category"
$wsDeposit.Range("I1").Value = "date"
$wsDeposit.Range("J1").Value = "legislator_name"
$wsDeposit.Range("K1").Value = "legislator_id"
$wsDeposit.Range("L1").Value = "source_file"
$wsDeposit.Range("M1").Value = "index"

# Make sure the shared-string cell F6 (formerly the text "1300000") becomes a
# real number, matching every other amount cell in the column.
$wsDeposit.Range("F6").Value = 1300000

$indices = @(50, 51, 52, 53, 54, 55, 56, 57)
for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $wsDeposit.Range("B$r`:F$r").Copy()
    $wsDeposit.Range("G$r`:M$r").PasteSpecial(-4122)

    $wsDeposit.Range("G$r").Value = "deposit"
    $wsDeposit.Range("H$r").Value = "normal"
    $dateSource.Copy()
    $wsDeposit.Range("I$r").PasteSpecial(-4163)
    $wsDeposit.Range("J$r").Value = "段宜康"
    $wsDeposit.Range("K$r").Value = 917
    $wsDeposit.Range("L$r").Value = "tmp25ce1"
    $wsDeposit.Range("M$r").Value = $indices[$i]
}
